# Generate Report for Handback
# Adds a new handback row (3d9034ca-d4b7-4303-981f-bee0da96184d) to the
# Overview, zh-cn and de-de sheets, mirroring the existing 71fb0c0c... row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add()

$wsOverview.Range("A3").Value = "3d9034ca-d4b7-4303-981f-bee0da96184d.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-09-05 10:33:07"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/628d847701d8d7b43b8c10dad0f57a3c32968c47/e2e/3d9034ca-d4b7-4303-981f-bee0da96184d.md", "", "", "e2e\3d9034ca-d4b7-4303-981f-bee0da96184d.md")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add()

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'True"
$wsZhCn.Range("G3").Value = "3d9034ca-d4b7-4303-981f-bee0da96184d.217a4da1669ecd9b3d537b200ab05b96a9ad5b11.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-05 10:32:51"
$wsZhCn.Range("J3").Value = "3d9034ca-d4b7-4303-981f-bee0da96184d.217a4da1669ecd9b3d537b200ab05b96a9ad5b11.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-05 10:34:00"
$wsZhCn.Range("L3").Value = "'"
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("N3").Value = "'"
$wsZhCn.Range("O3").Value = "'False"
$wsZhCn.Range("P3").Value = "'"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/628d847701d8d7b43b8c10dad0f57a3c32968c47/e2e/3d9034ca-d4b7-4303-981f-bee0da96184d.md", "", "", "3d9034ca-d4b7-4303-981f-bee0da96184d.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/17558788c995ce75727c261fbe0dfec57fc1a579/e2e/3d9034ca-d4b7-4303-981f-bee0da96184d.md", "", "", "3d9034ca-d4b7-4303-981f-bee0da96184d.md")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add()

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'True"
$wsDeDe.Range("G3").Value = "3d9034ca-d4b7-4303-981f-bee0da96184d.217a4da1669ecd9b3d537b200ab05b96a9ad5b11.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-05 10:33:07"
$wsDeDe.Range("J3").Value = "3d9034ca-d4b7-4303-981f-bee0da96184d.217a4da1669ecd9b3d537b200ab05b96a9ad5b11.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-05 10:34:27"
$wsDeDe.Range("L3").Value = "'"
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("N3").Value = "'"
$wsDeDe.Range("O3").Value = "'False"
$wsDeDe.Range("P3").Value = "'"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/628d847701d8d7b43b8c10dad0f57a3c32968c47/e2e/3d9034ca-d4b7-4303-981f-bee0da96184d.md", "", "", "3d9034ca-d4b7-4303-981f-bee0da96184d.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/fa94f5fb2ac079719c4db7aae33132cbad02e3a2/e2e/3d9034ca-d4b7-4303-981f-bee0da96184d.md", "", "", "3d9034ca-d4b7-4303-981f-bee0da96184d.md")
